$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting of the existing header cell (H1) for I1 and J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Values for column I (I0) and column J (IF), rows 2-36
$values = @{
    2  = @(6, 8)
    3  = @(7, 8)
    4  = @(6, 6)
    5  = @(5, 7)
    6  = @(5, 7)
    7  = @(5, 6)
    8  = @(3, 4)
    9  = @(6, 7)
    10 = @(6, 7)
    11 = @(7, 9)
    12 = @(6, 7)
    13 = @(2, 4)
    14 = @(9, 9)
    15 = @(8, 9)
    16 = @(7, 8)
    17 = @(2, 5)
    18 = @(2, 5)
    19 = @(2, 5)
    20 = @(1, 4)
    21 = @(1, 6)
    22 = @(1, 6)
    23 = @(1, 4)
    24 = @(1, 6)
    25 = @(1, 5)
    26 = @(1, 5)
    27 = @(1, 5)
    28 = @(1, 5)
    29 = @(1, 5)
    30 = @(1, 3)
    31 = @(1, 4)
    32 = @(1, 4)
    33 = @(1, 5)
    34 = @(1, 4)
    35 = @(1, 3)
    36 = @(1, 2)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
